{"js": "// Brand Statement edit \u2014 Office.js (Word JavaScript API)\n//\n// The author rewrote the single paragraph of the brand statement:\n//   - split the opening sentence (\"...Irvine.\") from the \"Quick learner...\"\n//     sentence into two runs, and reworded the certification sentence;\n//   - inserted a brand-new block (\"I enjoy the various opportunities...\n//     for both sellers and buyers.\") right after \"...various datasets. \";\n//   - inserted a new \"My adaptive ability...\" closing block right before\n//     \"Recently c...\", and removed the old closing block (which used to\n//     sit after the _GoBack bookmark), moving the bookmark to the very\n//     end of the paragraph.\n//\n// Because the whole paragraph's wording/run-layout changes, the most\n// faithful way to reproduce the target OOXML (matching run boundaries)\n// is to rebuild that single paragraph via insertOoxml(\"Replace\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = paragraphs.items[0];\n\nconst runsText = [\n  \"Data analyst with a certificate from the data analytics bootcamp from the University of California, Irvine. \",\n  \"Quick learner with the ability to apply knowledge of VBA, Python, SQL, JavaScript, R, and Tableau to analyze and visualize data from \",\n  \"various\",\n  \" datasets. \",\n  \"I enjoy\",\n  \" the various opportunities \",\n  \"provided\",\n  \" in this field, from storytelling to picture recognition. \",\n  \"My main \",\n  \"objective\",\n  \" \",\n  \"is to \",\n  \"work in advertising and determine the correct time to provide consumers with information for maximum \",\n  \"gain \",\n  \"for\",\n  \" both sellers and buyers. \",\n  \"Recently c\",\n  \"ompleted a project with a group that utilized machine learning methods to capture images of playing cards and \",\n  \"classify the number and suit.\",\n  \" \",\n  \"My adaptive ability to learn and \",\n  \"apply \",\n  \"new information and apply \",\n  \"them\",\n  \" \",\n  \"the highest level has set a foundation for me to work in exciting and fast-paced environments.\",\n];\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\")\n    .replace(/'/g, \"&apos;\");\n}\n\nconst runsXml = runsText\n  .map((t) => {\n    const preserve = /^\\s|\\s$|^$/.test(t) ? ' xml:space=\"preserve\"' : \"\";\n    return `<w:r><w:t${preserve}>${escapeXml(t)}</w:t></w:r>`;\n  })\n  .join(\"\");\n\nconst bookmarkXml =\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>';\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>${runsXml}${bookmarkXml}</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\ntargetParagraph.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Brand Statement edit \u2014 Word COM interop (PowerShell-style)\n#\n# The author rewrote the single paragraph of the brand statement:\n#   - split the opening sentence (\"...Irvine.\") from the \"Quick learner...\"\n#     sentence into two runs, and reworded the certification sentence;\n#   - inserted a brand-new block (\"I enjoy the various opportunities...\n#     for both sellers and buyers.\") right after \"...various datasets. \";\n#   - inserted a new \"My adaptive ability...\" closing block right before\n#     \"Recently c...\", and removed the old closing block (which used to\n#     sit after the _GoBack bookmark), moving the bookmark to the very\n#     end of the paragraph.\n#\n# Because the whole paragraph's wording/run-layout changes, the most\n# faithful way to reproduce the target OOXML (matching run boundaries) is\n# to rebuild that single paragraph's contents via Range.InsertXML \u2014 the\n# COM analogue of Office.js's Range.insertOoxml(\"Replace\").\n\n$d = $word.ActiveDocument\n\n$targetParagraph = $d.Paragraphs(1)\n\n# Range over the paragraph's contents, EXCLUDING the trailing paragraph\n# mark, so InsertXML replaces the runs in place and keeps this as the\n# one-and-only paragraph (instead of inserting a sibling paragraph).\n$r = $d.Range(0, $targetParagraph.Range.End - 1)\n\n$runsXml = \"\"\n$runsXml += '<w:r><w:t xml:space=\"preserve\">Data analyst with a certificate from the data analytics bootcamp from the University of California, Irvine. </w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\">Quick learner with the ability to apply knowledge of VBA, Python, SQL, JavaScript, R, and Tableau to analyze and visualize data from </w:t></w:r>'\n$runsXml += '<w:r><w:t>various</w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\"> datasets. </w:t></w:r>'\n$runsXml += '<w:r><w:t>I enjoy</w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\"> the various opportunities </w:t></w:r>'\n$runsXml += '<w:r><w:t>provided</w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\"> in this field, from storytelling to picture recognition. </w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\">My main </w:t></w:r>'\n$runsXml += '<w:r><w:t>objective</w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\">is to </w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\">work in advertising and determine the correct time to provide consumers with information for maximum </w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\">gain </w:t></w:r>'\n$runsXml += '<w:r><w:t>for</w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\"> both sellers and buyers. </w:t></w:r>'\n$runsXml += '<w:r><w:t>Recently c</w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\">ompleted a project with a group that utilized machine learning methods to capture images of playing cards and </w:t></w:r>'\n$runsXml += '<w:r><w:t>classify the number and suit.</w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\">My adaptive ability to learn and </w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\">apply </w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\">new information and apply </w:t></w:r>'\n$runsXml += '<w:r><w:t>them</w:t></w:r>'\n$runsXml += '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>'\n$runsXml += '<w:r><w:t>the highest level has set a foundation for me to work in exciting and fast-paced environments.</w:t></w:r>'\n\n$bookmarkXml = '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>'\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n  '<pkg:xmlData>' + `\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n  '<w:body>' + `\n  '<w:p>' + $runsXml + $bookmarkXml + '</w:p>' + `\n  '</w:body>' + `\n  '</w:document>' + `\n  '</pkg:xmlData>' + `\n  '</pkg:part>' + `\n  '</pkg:package>'\n\n$r.InsertXML($xml)\n"}
